$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-company"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet ---
$elem = $wb.Worksheets.Item("Elements")
$elem.Range("AI2").Value = ""
$elem.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-company"
